# Apply "repull data, push all data, mean calculation" updates.
# These edits change the dSF column (F) values for a number of rows
# on Sheet1, while leaving dS0 (E) and everything else untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$updates = @{
    3  = -4
    4  = -1
    6  = -2
    11 = 3
    12 = 1
    16 = 2
    21 = -1
    25 = -1
    26 = 3
    27 = 1
    30 = 2
    37 = 2
    42 = 1
    45 = -3
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
